$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (mirrors existing formatting of rows 2-4)
$data = @(
    @{ Row = 5; A = 10; B = "M931252509029" },
    @{ Row = 6; A = 4;  B = "N304350709089" },
    @{ Row = 7; A = 6;  B = "R931101109037" }
)

foreach ($entry in $data) {
    $r = $entry.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $entry.A

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $entry.B

    # Copy the formatting from row 4 (last existing data row) onto the new row,
    # so style, borders, alignment stay consistent with existing data rows.
    $ws.Range("A4").Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$excel.CutCopyMode = $false
